$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(178, "Monday, Jan 16", "5:40 AM", "LO3942", "Warsaw",       "(WAW)", "LOT ",      "E190", "(SP-LMH)", "5:58 AM", "0 hours, 18 minutes"),
    @(179, "Monday, Jan 16", "6:20 AM", "FR3594", "Milan",        "(BGY)", "Ryanair ",  "B738", "(SP-RSM)", "6:29 AM", "0 hours, 9 minutes"),
    @(180, "Monday, Jan 16", "6:30 AM", "FR1751", "London",       "(STN)", "Ryanair ",  "B738", "(SP-RKR)", "6:45 AM", "0 hours, 15 minutes"),
    @(181, "Monday, Jan 16", "6:30 AM", "FR7938", "Edinburgh",    "(EDI)", "Ryanair ",  "B738", "(SP-RSX)", "6:42 AM", "0 hours, 12 minutes"),
    @(182, "Monday, Jan 16", "6:35 AM", "LH1381", "Frankfurt",    "(FRA)", "Lufthansa ","CRJ9", "(D-ACNE)", "7:06 AM", "0 hours, 31 minutes"),
    @(183, "Monday, Jan 16", "8:05 AM", "LH1641", "Munich",       "(MUC)", "Lufthansa ","CRJ9", "(D-ACNM)", "8:25 AM", "0 hours, 20 minutes"),
    @(184, "Monday, Jan 16", "9:15 AM", "UNKNOWN","Fuerteventura","(FUE)", "Enter Air ","B738", "(SP-ENL)", "9:21 AM", "0 hours, 6 minutes")
)

$startRow = 179
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
    $ws.Cells.Item($r, 6).Value = $values[5]
    $ws.Cells.Item($r, 7).Value = $values[6]
    $ws.Cells.Item($r, 8).Value = $values[7]
    $ws.Cells.Item($r, 9).Value = $values[8]
    $ws.Cells.Item($r, 10).Value = $values[9]
    $ws.Cells.Item($r, 12).Value = $values[10]
}
